$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Executive Summary
# ---------------------------------------------------------------------------
$wsExec = $wb.Worksheets.Item("Executive Summary")
$wsExec.Range("B2").Value = 160
$wsExec.Range("B3").Value = 292.5
$wsExec.Range("C3").Value = 82.8
$wsExec.Range("B4").Value = 468.8
$wsExec.Range("C4").Value = 60.3
$wsExec.Range("B5").Value = 663.8
$wsExec.Range("C5").Value = 41.6
$wsExec.Range("B6").Value = 837.2
$wsExec.Range("C6").Value = 26.1
$wsExec.Range("B7").Value = 964.5
$wsExec.Range("C7").Value = 15.2

# ---------------------------------------------------------------------------
# Sheet: Market Calculations - Generation Z rows (2025-2030)
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Market Calculations")

# F-column cells hold "$940"-style text; Excel auto-converts pure currency-
# looking strings to numbers on assignment, so force text storage via a
# temporary Text number format, then restore the default (Normal) style so
# the cell ends up exactly like the original - plain text, no explicit style.
$genZSpendRows = @(6, 12, 18, 24, 30, 36)
foreach ($r in $genZSpendRows) {
    $cell = $wsCalc.Cells.Item($r, 6)
    $cell.NumberFormat = "@"
    $cell.Value = "$3,450"
    $cell.Style = "Normal"
}

$wsCalc.Range("G6").Value = 8.737
$wsCalc.Range("H6").Value = "67.0M × 21% × 18% × $3,450"

$wsCalc.Range("G12").Value = 16.99
$wsCalc.Range("H12").Value = "67.0M × 21% × 35% × $3,450"

$wsCalc.Range("G18").Value = 26.698
$wsCalc.Range("H18").Value = "67.0M × 21% × 55% × $3,450"

$wsCalc.Range("G24").Value = 34.95
$wsCalc.Range("H24").Value = "67.0M × 21% × 72% × $3,450"

$wsCalc.Range("G30").Value = 39.804
$wsCalc.Range("H30").Value = "67.0M × 21% × 82% × $3,450"

$wsCalc.Range("G36").Value = 42.231
$wsCalc.Range("H36").Value = "67.0M × 21% × 87% × $3,450"

# ---------------------------------------------------------------------------
# Sheet: Generation Breakdown
# ---------------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("Generation Breakdown")

$wsGen.Range("B6").Value = 8.74
$wsGen.Range("C6").Value = 16.99
$wsGen.Range("D6").Value = 26.7
$wsGen.Range("E6").Value = 34.95
$wsGen.Range("F6").Value = 39.8
$wsGen.Range("G6").Value = 42.23

$wsGen.Range("B8").Value = 160
$wsGen.Range("C8").Value = 292.5
$wsGen.Range("D8").Value = 468.8
$wsGen.Range("E8").Value = 663.8
$wsGen.Range("F8").Value = 837.2
$wsGen.Range("G8").Value = 964.5

# ---------------------------------------------------------------------------
# Sheet: Input Parameters
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input Parameters")

$cellC18 = $wsInput.Range("C18")
$cellC18.NumberFormat = "@"
$cellC18.Value = "$3,450"
$cellC18.Style = "Normal"

$wsInput.Range("D18").Value = "Piper Sandler/NIQ/DataReportal"
$wsInput.Range("E18").Value = "https://www.pipersandler.com/teens"

# ---------------------------------------------------------------------------
# Sheet: Assumptions
# ---------------------------------------------------------------------------
$wsAssum = $wb.Worksheets.Item("Assumptions")

$wsAssum.Range("B7").Value = "$3,450 based on total online spending (80% of shopping)"
$wsAssum.Range("C7").Value = "Calculated from Piper Sandler teen data, NIQ online shopping patterns, and DataReportal averages"

# ---------------------------------------------------------------------------
# Sheet: Data Sources
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data Sources")

$wsData.Range("A9").Value = "Gen Z Total Online Spending"
$wsData.Range("B9").Value = "Piper Sandler/NIQ/DataReportal"
$wsData.Range("C9").Value = "https://www.pipersandler.com/teens"
$wsData.Range("D9").Value = "2024-2025"
